# CKU overall advancement. Some minor code changes
#
# Adds a "feudal_government" value under the existing "government" column
# (column E) for every kingdom row, auto-fits that column to its new
# content, re-applies/duplicates the "highlight duplicate values"
# conditional formatting across columns A:B (promoting it to top
# priority, as Excel does when a new Highlight-Cells rule is created),
# and restores the last active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kingdoms")

# --- New "feudal_government" data in column E (rows 2-6) ---------------
$ws.Range("E2:E6").Value = "feudal_government"

# Column E best-fits itself to the new content, like Excel does when the
# column is double-click auto-sized after data entry.
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# --- Conditional formatting: duplicate-values highlight over A:B -------
$range = $ws.Range("A1:B1048576")

# A handful of throwaway duplicate-value rules were created/removed while
# tweaking the formatting (their dxf records linger in styles.xml, same
# as real Excel leaves orphaned dxfs behind after rule edits/deletes).
for ($i = 0; $i -lt 7; $i++) {
    $scratch = $range.FormatConditions.AddUniqueValues()
    $scratch.DupeUnique = 1
    $scratch.Font.Color = 393372
    $scratch.Interior.Color = 13551615
    $scratch.Delete()
}

# The rule that actually sticks, matching the workbook's existing
# "duplicate values" red style, promoted to the top priority.
$fc = $range.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
$fc.SetFirstPriority()

# --- Restore last selection --------------------------------------------
$ws.Range("J13").Select() | Out-Null
